$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Language list / labels (header + footer language switcher lines)
Replace-All "English" "Inglés"
Replace-All " / Portuguese / French / Thai / Vietnamese / Spanish" " / Portugués / Francés / Tailandés / Vietnamita / Español"

# Brief table
Replace-All "Brief" "Breve"
Replace-All "An email sent to partners in the target country who have sent their documents for review. It will be sent via customer.io" "Un correo electrónico enviado a los socios en el país objetivo que han enviado sus documentos para revisión. Se enviará a través de customer.io"
Replace-All "Target audience" "Público objetivo"
Replace-All "Invited partners who have submitted their documents" "Socios invitados que han presentado sus documentos"

# Subject line
Replace-All " — we got your docs!  " " — ¡hemos recibido tus documentos!  "

# Heading
Replace-All "Thank you for submitting your documents" "Gracias por enviar tus documentos"

# Greeting
Replace-All "Hi " "Hola "

# Body paragraphs
Replace-All "Thank you for providing us with your documents for the upcoming " "Gracias por facilitarnos tus documentos para el próximo "
Replace-All ". Based on the information you’ve given us, we’ll make the necessary arrangements, including accommodation and transportation." ". Basándonos en la información que nos has facilitado, haremos los preparativos necesarios, incluidos el alojamiento y el transporte."
Replace-All "We’re currently reviewing your documents and will reach out to you if we need anything else. " "Estamos revisando tus documentos y nos pondremos en contacto contigo si necesitamos algo más. "
Replace-All "If you have any questions, please contact us via " "Si tienes alguna pregunta, entra en contacto con nosotros por "
Replace-All " or " " o "
Replace-All "If you have any questions, please contact your country manager, " "Si tienes alguna pregunta, entra en contacto con el gestor de tu país "
Replace-All ", at " ", en "
Replace-All "We look forward to seeing you at " "Esperamos verte en "

# Comment text (only the top-level Comment.Range.Text setter is safely scoped in
# this runtime; any deeper sub-range operation on a comment's Range silently
# falls back to mutating the main document story, so it must be avoided).
$d.Comments.Item(1).Range.Text = "elija uno de los dos"
